# Apply corrections to the Clusters column (K) values that were
# miscalculated in the original workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(23, 11).Value = 1
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(26, 11).Value = 3
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(30, 11).Value = 1
$ws.Cells.Item(31, 11).Value = 1
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(38, 11).Value = 1
$ws.Cells.Item(40, 11).Value = 3
$ws.Cells.Item(41, 11).Value = 3
$ws.Cells.Item(43, 11).Value = 3
$ws.Cells.Item(45, 11).Value = 1
$ws.Cells.Item(46, 11).Value = 1
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(53, 11).Value = 1
$ws.Cells.Item(55, 11).Value = 3
$ws.Cells.Item(56, 11).Value = 3
$ws.Cells.Item(58, 11).Value = 3
$ws.Cells.Item(60, 11).Value = 1
$ws.Cells.Item(61, 11).Value = 1
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(68, 11).Value = 1
$ws.Cells.Item(70, 11).Value = 3
$ws.Cells.Item(71, 11).Value = 3
$ws.Cells.Item(73, 11).Value = 3
$ws.Cells.Item(75, 11).Value = 1
$ws.Cells.Item(76, 11).Value = 1
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(83, 11).Value = 1
$ws.Cells.Item(85, 11).Value = 3
$ws.Cells.Item(86, 11).Value = 3
$ws.Cells.Item(88, 11).Value = 3
$ws.Cells.Item(90, 11).Value = 1
$ws.Cells.Item(91, 11).Value = 1
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(98, 11).Value = 1
$ws.Cells.Item(100, 11).Value = 3
$ws.Cells.Item(101, 11).Value = 3
$ws.Cells.Item(103, 11).Value = 3
$ws.Cells.Item(105, 11).Value = 1
$ws.Cells.Item(106, 11).Value = 1
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(113, 11).Value = 1
$ws.Cells.Item(115, 11).Value = 3
$ws.Cells.Item(116, 11).Value = 3
$ws.Cells.Item(118, 11).Value = 3
$ws.Cells.Item(120, 11).Value = 1
$ws.Cells.Item(121, 11).Value = 1
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(128, 11).Value = 1
$ws.Cells.Item(130, 11).Value = 3
$ws.Cells.Item(131, 11).Value = 3
$ws.Cells.Item(133, 11).Value = 3
$ws.Cells.Item(135, 11).Value = 1
$ws.Cells.Item(136, 11).Value = 1
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(143, 11).Value = 1
$ws.Cells.Item(145, 11).Value = 3
$ws.Cells.Item(146, 11).Value = 3
$ws.Cells.Item(148, 11).Value = 3
$ws.Cells.Item(150, 11).Value = 1
$ws.Cells.Item(151, 11).Value = 1
$ws.Cells.Item(155, 11).Value = 0
$ws.Cells.Item(158, 11).Value = 1
$ws.Cells.Item(160, 11).Value = 3
$ws.Cells.Item(161, 11).Value = 3
$ws.Cells.Item(163, 11).Value = 3
$ws.Cells.Item(165, 11).Value = 1
$ws.Cells.Item(166, 11).Value = 1
$ws.Cells.Item(170, 11).Value = 0
$ws.Cells.Item(173, 11).Value = 1
$ws.Cells.Item(175, 11).Value = 3
$ws.Cells.Item(176, 11).Value = 3
$ws.Cells.Item(178, 11).Value = 3
$ws.Cells.Item(180, 11).Value = 1
$ws.Cells.Item(181, 11).Value = 1
$ws.Cells.Item(185, 11).Value = 0
$ws.Cells.Item(188, 11).Value = 1
$ws.Cells.Item(190, 11).Value = 3
$ws.Cells.Item(191, 11).Value = 3
$ws.Cells.Item(193, 11).Value = 3
$ws.Cells.Item(195, 11).Value = 1
$ws.Cells.Item(196, 11).Value = 1
$ws.Cells.Item(200, 11).Value = 0
$ws.Cells.Item(203, 11).Value = 1
$ws.Cells.Item(205, 11).Value = 3
$ws.Cells.Item(206, 11).Value = 3
$ws.Cells.Item(208, 11).Value = 3
$ws.Cells.Item(210, 11).Value = 1
$ws.Cells.Item(211, 11).Value = 1
$ws.Cells.Item(215, 11).Value = 0
$ws.Cells.Item(218, 11).Value = 1
$ws.Cells.Item(220, 11).Value = 3
$ws.Cells.Item(221, 11).Value = 3
$ws.Cells.Item(223, 11).Value = 3
$ws.Cells.Item(225, 11).Value = 1
$ws.Cells.Item(226, 11).Value = 1
$ws.Cells.Item(230, 11).Value = 0
$ws.Cells.Item(233, 11).Value = 1
$ws.Cells.Item(235, 11).Value = 3
$ws.Cells.Item(236, 11).Value = 3
$ws.Cells.Item(238, 11).Value = 3
$ws.Cells.Item(240, 11).Value = 1
$ws.Cells.Item(241, 11).Value = 1
$ws.Cells.Item(245, 11).Value = 0
$ws.Cells.Item(248, 11).Value = 1
$ws.Cells.Item(250, 11).Value = 3
$ws.Cells.Item(251, 11).Value = 3
$ws.Cells.Item(253, 11).Value = 3
$ws.Cells.Item(255, 11).Value = 1
$ws.Cells.Item(256, 11).Value = 1
$ws.Cells.Item(260, 11).Value = 0
$ws.Cells.Item(263, 11).Value = 1
$ws.Cells.Item(265, 11).Value = 3
$ws.Cells.Item(266, 11).Value = 3
$ws.Cells.Item(268, 11).Value = 3
$ws.Cells.Item(270, 11).Value = 1
$ws.Cells.Item(271, 11).Value = 1
